# Correct status name labels across the worksheet.
# Replaces the old "statut_label"/"statut_name" text values with the corrected wording,
# matching the updated shared-string text from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old exact cell text -> new exact cell text.
$replacements = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$firstRow = $used.Row
$firstCol = $used.Column

for ($r = $firstRow; $r -lt ($firstRow + $rowCount); $r++) {
    for ($c = $firstCol; $c -lt ($firstCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value2 = $replacements[$val]
        }
    }
}
